# Update portfolio presentation: change project section on slide 3
# from "그로스폴리오 캠페인" (marketing campaign project) to
# "오류수정하기 프로젝트" (bug-fix project), updating the title,
# description bullets, and achievement bullets while preserving all
# paragraph-level formatting (font size, bold, spacing, etc.).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

$tr.Paragraphs(1,1).Text  = "오류수정하기 프로젝트"
$tr.Paragraphs(3,1).Text  = "• 서비스 내 사용자 경험을 저해하는 버그 및 오류 사항 모니터링 및 수집"
$tr.Paragraphs(4,1).Text  = "• 고객 피드백 데이터 분석을 통한 우선순위 오류 항목 선별"
$tr.Paragraphs(5,1).Text  = "• 개발팀과의 협업을 통한 오류 수정 프로세스 수립 및 관리"
$tr.Paragraphs(6,1).Text  = "• 수정된 오류에 대한 사용자 만족도 조사 및 피드백 수집"
$tr.Paragraphs(7,1).Text  = "• 주간 오류 리포트 작성 및 개선사항 도출"
$tr.Paragraphs(9,1).Text  = "• 서비스 오류 신고 건수 전월 대비 45% 감소"
$tr.Paragraphs(10,1).Text = "• 사용자 서비스 만족도 조사 결과 85점으로 향상"
$tr.Paragraphs(11,1).Text = "• 플랫폼 안정성 지표 92%까지 개선"
